$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the "Animation" row (old row 47), to hold the
# new "Cone" / "Cylindre" sub-items of the Object3D group (Cube, Sphere, ...).
$ws.Rows("47:48").Insert()

# Populate the two newly inserted rows (A47/A48 already inherited style "2"
# from the row above as part of the insert).
$ws.Range("B47").Value = "Cone"
$ws.Range("B48").Value = "Cylindre"

# Extend the Object3D group merge (A45:A46) down to cover the new rows.
$ws.Range("A45:A48").Merge()

# Move the "Yes" status value that used to sit next to "Tilemap" (old C48,
# now shifted to C50) up to sit next to "Animation" (now row 49) instead.
$ws.Range("C50").Value = ""
$ws.Range("C49").Value = "Yes"

# Update the active selection to match the saved view.
$ws.Range("B48").Select()
